# autodiagnostico.xlsx adjustments
# - Fills in missing "Pasos a Seguir" / "Datos de Prueba" steps for DMZ, IPv4 Port
#   Mapping, Reserva DHCP and Dispositivos Conectados test cases (rows 8-11)
# - Fills in "Resultado Esperado" / "Resultado Obtenido" for the order-creation
#   test case (row 6)
# - Adjusts row heights that grew because of the new, longer wrapped text
# - Updates the active selection left on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (CP_AUTO_007 / funcion DMZ) ---------------------------------
$ws.Range("F8").Value = "1.Clic en Botón Opciones`n2.Clic en Opción DMZ`n3.Clic casilla Habilitar DMZ`n4.Diligenciar IP aleatoria`n5.Clic en Botón Refrescar`n6.Clic en Botón Cancelar`n7."
$ws.Range("G8").Value = "N/A"
$ws.Range("G8").Font.Underline = $true
$ws.Rows.Item(8).RowHeight = 58.5

# --- Row 9 (CP_AUTO_008 / ipv4 port Mapping) ---------------------------
$ws.Range("F9").Value = "1.Botón ""Opciones""`n2.Opción ""IPv4 Port Mapping""`n3.Clic en el campo ""Protocolo""`n4.Seleccionar opción aleatoria en ""Protocolo""`n5.Diligenciar campo ""Dirección IP"" con una IPv4 aleatoria`n6.Clic en botón ""Refrescar""`n7.Clic en botón ""Cancelar"""
$ws.Range("G9").Value = "N/A"
$ws.Range("G9").Font.Underline = $true
$ws.Rows.Item(9).RowHeight = 142.5

# --- Row 10 (CP_AUTO_009 / reserva DHCP) -------------------------------
$ws.Range("F10").Value = "1.Clic en botón ""Opciones""`n2.Clic en Opción ""Reserva DHCP""`n3.Diligenciar MAC aleatoria`n4.Diligenciar IPv4 aleatoria`n5.Clic en botón ""Refrescar""`n6.Clic en botón ""Cancelar"""
$ws.Range("G10").Value = "N/A"
$ws.Range("G10").Font.Underline = $true
$ws.Rows.Item(10).RowHeight = 85.5

# --- Row 11 (CP_AUTO_010 / Dispositivos conectados) --------------------
$ws.Range("F11").Value = "1.Clic en botón ""Opciones""`n2.Opción ""Dispositivos Conectados""`n3.Clic en flecha desplegable del primer dispositivo`n4.Clic en flecha desplegable del segundo dispositivo`n5.Clic en botón ""Recargar/Refrescar""`n6.Cerrar modal ""Dispositivos Conectados"""
$ws.Range("G11").Value = "N/A"
$ws.Range("G11").Font.Underline = $true
$ws.Rows.Item(11).RowHeight = 134.25

# --- Row 6 (CP_AUTO_006 / funcion UPnP) --------------------------------
$ws.Range("H6").Value = "El sistema debe permitir la creacion de la orden"
$ws.Range("I6").Value = "La orden se crea correctamente"

# --- Leave the selection where the author left it in the saved file ----
$ws.Range("J6").Select()
